$d = $word.ActiveDocument

# 1. Remove the whole paragraph "FireMagic(TM) 3D irrigado + equipo de irrigação." (including its paragraph mark)
$d.Content.Find.Execute("FireMagic" + [char]0x2122 + " 3D irrigado + equipo de irrigação.`r", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Prefix the material list items with a bullet character "• "
$items = @(
    "Cateter de Ablação Irrigado",
    "Equipo de Irrigação",
    "Cateter Decapolar",
    "Cateter Quadripolar",
    "Introdutor"
)

foreach ($item in $items) {
    $d.Content.Find.Execute($item, $true, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " " + $item, 2)
}
